$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price column cells being updated to remain plain text (matches
# the original inline-string cell type) instead of being auto-parsed as numbers.
$ws.Range("D2:D23").NumberFormat = "@"
$ws.Range("D25:D38").NumberFormat = "@"
$ws.Range("D40:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.210.47"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").Value = "1.916.72"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -1.30%  "
$ws.Range("D5").Value = "327.99"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "0.4680"
$ws.Range("E7").Value = "  -5.98%  "
$ws.Range("D8").Value = "0.4012"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").Value = "52.84"
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("D10").Value = "0.08419"
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").Value = "1.048"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("D12").Value = "22.18"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "1.898.05"
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").Value = "7.449"
$ws.Range("E14").Value = "  -6.83%  "
$ws.Range("D15").Value = "6.079"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "89.68"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "0.00001069"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").Value = "0.06608"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "18.00"
$ws.Range("E20").Value = "  -7.76%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "5.749"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").Value = "28.166.74"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("E24").Value = "  -6.54%  "
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "2.113.95"
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").Value = "153.47"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").Value = "5.772"
$ws.Range("E29").Value = "  -8.51%  "
$ws.Range("D30").Value = "2.135"
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").Value = "123.60"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").Value = "0.9788"
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("D33").Value = "0.09673"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "1.443"
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("D35").Value = "3.647"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "5.554"
$ws.Range("E36").Value = "  -4.66%  "
$ws.Range("D37").Value = "8.839"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").Value = "1.266"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("D40").Value = "0.06186"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("D41").Value = "0.6171"
$ws.Range("E41").Value = "  -4.86%  "
$ws.Range("D42").Value = "11.06"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "0.1907"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").Value = "1.302"
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("D46").Value = "0.5862"
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("D47").Value = "12.80"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("D48").Value = "2.030"
$ws.Range("E48").Value = "  -6.63%  "
$ws.Range("D49").Value = "3.438"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").Value = "0.06903"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "111.21"
$ws.Range("E51").Value = "  -1.60%  "
